$d = $word.ActiveDocument

function Find-ParaIndex($pattern) {
    $idx = 1
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like $pattern) {
            return $idx
        }
        $idx = $idx + 1
    }
    return -1
}

# ---------------------------------------------------------------------------
# Change 1: the empty paragraph right before "Icons need to be added..."
# merges with that paragraph (the paragraph break between them is removed);
# the tail of the old sentence is then replaced by a new sentence about
# company/education logos, kept as its own run; finally a brand-new
# paragraph with a further sentence about "My own name..." is appended
# right after it.
# ---------------------------------------------------------------------------

$iconsIdx = Find-ParaIndex "Icons need to be added*"
$prevPara = $d.Paragraphs.Item($iconsIdx - 1)
$prevPara.Range.Delete()

# The empty paragraph above is now gone, so "Icons..." shifted up by one.
$iconsIdx = Find-ParaIndex "Icons need to be added*"
$iconsPara = $d.Paragraphs.Item($iconsIdx)

# Drop the old tail of the sentence, keeping the shared lead-in text.
$iconsPara.Range.Find.Execute("needs to look good in it’s current setting.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

$iconsPara = $d.Paragraphs.Item($iconsIdx)

# Insert the new sentence as its own run. Toggling Bold on and back off
# keeps Word from silently coalescing it into the previous run even though
# the two runs end up with visually identical formatting.
$secondSentence = "I want the companies and educational systems to have logo’s of them to the left side of the orange column"
$insertStart = $iconsPara.Range.End - 1
$insertPoint = $d.Range($insertStart, $insertStart)
$insertPoint.InsertAfter($secondSentence)

$iconsPara = $d.Paragraphs.Item($iconsIdx)
$secondRunRange = $d.Range($insertStart, $iconsPara.Range.End - 1)
$secondRunRange.Font.Bold = 1
$secondRunRange.Font.Bold = 0

# Add a new paragraph right after with the additional sentence about the name.
$iconsPara = $d.Paragraphs.Item($iconsIdx)
$iconsPara.Range.InsertParagraphAfter()
$namePara = $d.Paragraphs.Item($iconsIdx + 1)
$namePara.Range.Text = "My own name should be in the correct orange font, same as last name"

# ---------------------------------------------------------------------------
# Change 2: a brand new paragraph right after the "sticky sessions" one.
# ---------------------------------------------------------------------------
$stickyIdx = Find-ParaIndex "It*ll need sticky sessions*"
$stickyPara = $d.Paragraphs.Item($stickyIdx)
$stickyPara.Range.InsertParagraphAfter()
$cleanupPara = $d.Paragraphs.Item($stickyIdx + 1)
$cleanupPara.Range.Text = "Cleanup github repo to conly contain 1 copmiler project domein name"

Write-Output "OK"
